$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.600.98"
$ws.Range("E2").Value = "  -0.31%  "
$ws.Range("D3").Value = "2.277.10"
$ws.Range("E3").Value = "  -0.40%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "122.80"
$ws.Range("E5").Value = "  +6.25%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "266.38"
$ws.Range("E6").Value = "  -0.75%  "
$ws.Range("E7").Value = "  +2.60%  "
$ws.Range("E8").Value = "  +0.22%  "
$ws.Range("E9").Value = "  +0.81%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "47.96"
$ws.Range("E10").Value = "  -2.40%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0946"
$ws.Range("E11").Value = "  +0.36%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "9.26"
$ws.Range("E12").Value = "  +3.81%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.106"
$ws.Range("E13").Value = "  -1.20%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "15.45"
$ws.Range("E14").Value = "  -2.72%  "
$ws.Range("E15").Value = "  +3.20%  "
$ws.Range("D16").Value = "2.618.10"
$ws.Range("E16").Value = "  -0.62%  "
$ws.Range("D17").Value = "2.284.32"
$ws.Range("E17").Value = "  -0.63%  "
$ws.Range("D18").Value = "43.557.52"
$ws.Range("E18").Value = "  -0.14%  "
$ws.Range("E19").Value = "  +1.24%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.26"
$ws.Range("E21").Value = "  -0.04%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.40"
$ws.Range("E22").Value = "  -0.51%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "235.32"
$ws.Range("E23").Value = "  +1.07%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.56"
$ws.Range("E24").Value = "  -3.53%  "
$ws.Range("E25").Value = "  -1.15%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.98"
$ws.Range("E26").Value = "  +2.98%  "
$ws.Range("E27").Value = "  +1.70%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "42.18"
$ws.Range("E28").Value = "  +0.59%  "
$ws.Range("E29").Value = "  -0.74%  "
$ws.Range("E30").Value = "  +0.65%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "172.27"
$ws.Range("E31").Value = "  -0.46%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "21.71"
$ws.Range("E32").Value = "  +0.55%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0918"
$ws.Range("E34").Value = "  +0.45%  "
$ws.Range("B35").Value = "Stellar"
$ws.Range("C35").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.130"
$ws.Range("E35").Value = "  +2.05%  "
$ws.Range("B36").Value = "NEARProtocol"
$ws.Range("C36").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.31"
$ws.Range("E36").Value = "  +13.79%  "
$ws.Range("E37").Value = "  +4.93%  "
$ws.Range("E38").Value = "  -1.41%  "
$ws.Range("E39").Value = "  -1.17%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.55"
$ws.Range("E40").Value = "  +4.76%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "13.97"
$ws.Range("E41").Value = "  -4.26%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "74.07"
$ws.Range("E42").Value = "  -0.37%  "
$ws.Range("E43").Value = "  -1.32%  "
$ws.Range("E44").Value = "  -0.12%  "
$ws.Range("E45").Value = "  -0.20%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.72"
$ws.Range("E46").Value = "  -10.25%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "74.02"
$ws.Range("E47").Value = "  +40.83%  "
$ws.Range("E48").Value = "  -0.52%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.54"
$ws.Range("E49").Value = "  -1.68%  "
$ws.Range("E50").Value = "  +0.23%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "101.66"
$ws.Range("E51").Value = "  -1.01%  "
